$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BatteryPack-AZA")
$ws.Range("Q17").Value = "Test"
$ws.Range("Q17").Font.Bold = $true
